$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("example2")

# --- Update the three period-header cells (E1, G1, H1) on the second sheet ---
# These are plain text labels (not real dates), so we briefly force a text
# number format while writing the new value to stop the engine from
# auto-converting strings like "jan 2011" into date serial numbers.

# E1: "2010Q2" -> "apr 2010", and this one additionally gets a real date
# display format (YYYY-MM-DD) applied to it afterwards.
$ws2.Range("E1").NumberFormat = "@"
$ws2.Range("E1").Value = "apr 2010"

# G1: "2011Q1" -> "jan 2011" (keeps the default/general format)
$ws2.Range("G1").NumberFormat = "@"
$ws2.Range("G1").Value = "jan 2011"
$ws2.Range("G1").NumberFormat = "General"

# H1: "2011Q2" -> "apr 2011" (keeps the default/general format)
$ws2.Range("H1").NumberFormat = "@"
$ws2.Range("H1").Value = "apr 2011"
$ws2.Range("H1").NumberFormat = "General"

# Now apply the custom date format to E1 (adds numFmtId 165 "YYYY-MM-DD"
# to the stylesheet) while leaving its value as the literal text "apr 2010".
$ws2.Range("E1").NumberFormat = "YYYY\-MM\-DD"

# --- Update the selection / active cell on the second (active) sheet ---
$ws2.Range("H2").Select()
